$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1: "Save" - copy formatting from the last existing header
# cell (G1, bold/centered/bordered) so H1 reuses the same style index,
# then set its text.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Data cell H2: numeric 0, matching the plain (unstyled) data cells.
$ws.Range("H2").Value = 0
